$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset previously non-zero values to 0
$ws.Range("H2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("B49").Value = 0
$ws.Range("E65").Value = 0

# Set new non-zero values
$ws.Range("G62").Value = 0.03100000000000003
$ws.Range("B68").Value = 0.01100000000000001
$ws.Range("K68").Value = -0.01800000000000002
$ws.Range("F81").Value = 0.1019999999999999
$ws.Range("G81").Value = 0.1699999999999999
$ws.Range("H84").Value = 0.116
$ws.Range("C88").Value = -0.06799999999999995
$ws.Range("L89").Value = -0.01199999999999996
$ws.Range("F93").Value = -0.02600000000000002
$ws.Range("G93").Value = 0.008000000000000007
$ws.Range("H97").Value = 0.05299999999999999
$ws.Range("J97").Value = 0.006000000000000005
$ws.Range("L100").Value = -0.03299999999999997
